$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each affected cell is assigned with a leading apostrophe so Excel treats
# numeric-looking text (e.g. "313.01") as literal text instead of coercing it
# to a Double (which would corrupt values like "20.00" -> 20 or introduce
# floating-point noise). The Style reset afterward clears the quote-prefix
# formatting flag that this trick leaves behind, matching the unstyled source cells.

$ws.Range("D2").Value = "'27.189.68"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "'  +0.64%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Range("D3").Value = "'1.850.27"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "'  +1.12%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Range("E4").Value = "'  -0.50%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Range("D5").Value = "'313.01"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "'  +0.18%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Range("E6").Value = "'  -0.41%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Range("D8").Value = "'0.3705"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Range("D9").Value = "'0.07274"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = "'  -0.98%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Range("D10").Value = "'0.8830"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = "'  +0.89%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Range("D11").Value = "'20.00"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = "'  +0.87%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Range("D12").Value = "'0.07818"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "'  -1.60%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Range("D13").Value = "'1.927.72"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = "'  +5.29%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Range("D14").Value = "'5.369"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "'  +0.45%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Range("D15").Value = "'6.493"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "'  -0.95%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Range("D16").Value = "'91.24"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range("E16").Value = "'  -0.28%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Range("E17").Value = "'  -0.48%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Range("D18").Value = "'0.000008923"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "'  +0.19%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Range("E19").Value = "'  -0.29%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Range("D20").Value = "'14.68"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "'  -0.75%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Range("D21").Value = "'27.220.34"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "'  +0.55%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Range("D22").Value = "'5.054"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Range("E22").Value = "'  -1.23%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Range("E23").Value = "'  -0.77%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Range("D24").Value = "'2.126.03"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "'  +2.53%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Range("D25").Value = "'1.950"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "'  +5.55%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Range("D26").Value = "'151.59"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = "'  -1.14%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Range("D27").Value = "'18.37"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = "'  -0.37%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Range("D28").Value = "'2.063"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "'  +0.92%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Range("D29").Value = "'115.52"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = "'  +0.02%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Range("D30").Value = "'5.046"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("D31").Value = "'0.08822"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = "'  -0.96%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Range("D32").Value = "'3.094"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = "'  +4.36%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Range("D33").Value = "'0.7611"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range("E33").Value = "'  +3.85%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Range("D34").Value = "'1.166"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "'  +3.38%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Range("D35").Value = "'4.495"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "'  +1.36%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Range("D36").Value = "'2.725"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "'  +9.88%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Range("D37").Value = "'1.081"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("D38").Value = "'0.01940"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = "'  -0.67%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Range("D39").Value = "'0.05230"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "'  -0.05%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Range("D40").Value = "'2.944"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = "'  -0.09%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Range("D41").Value = "'7.064"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "'  -0.66%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Range("D42").Value = "'0.5089"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "'  -1.43%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Range("E43").Value = "'  -0.23%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Range("D44").Value = "'8.369"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = "'  +2.10%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Range("D45").Value = "'0.4778"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "'  -1.43%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Range("D46").Value = "'10.26"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "'  +0.88%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Range("E47").Value = "'  -0.47%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Range("E48").Value = "'  +0.16%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Range("D49").Value = "'1.633"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Range("E50").Value = "'  +0.19%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Range("D51").Value = "'65.66"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range("E51").Value = "'  +1.10%  "
$ws.Cells.Item(51, 5).Style = "Normal"
